$d = $word.ActiveDocument
$dash = [char]0x2013

# --- Change 1: paragraph 1 text -------------------------------------------------
# "This is a Microsoft word document." ->
#   "This is a Microsoft word document.  " (plain run, two trailing spaces)
#   + three red-colored runs spelling "(This is a change \u2013 Version for main branch)"
$p1 = $d.Paragraphs.Item(1)
$r1 = $p1.Range

# two trailing spaces, still default (black) formatting
$r1.InsertAfter("  ")

$seg1 = "(This is a change " + $dash + " Ve"
$seg2 = "rsion for main branch"
$seg3 = ")"

$insStart = $r1.End - 1
$r1.InsertAfter($seg1)
$run1 = $d.Range($insStart, $r1.End - 1)
$run1.Font.Color = 255

$insStart = $r1.End - 1
$r1.InsertAfter($seg2)
$run2 = $d.Range($insStart, $r1.End - 1)
$run2.Font.Color = 255

$insStart = $r1.End - 1
$r1.InsertAfter($seg3)
$run3 = $d.Range($insStart, $r1.End - 1)
$run3.Font.Color = 255

# --- Change 2: append a new, empty, shaded paragraph at the end of the body ----
$lastPara = $d.Paragraphs.Last
$lastRange = $lastPara.Range
$lastRange.Find.Execute("we are free at last.", $false, $false, $false, $false, $false, $true, 1, $false, "we are free at last.^p", 2)

$newPara = $d.Paragraphs.Last
$newPara.Style = "Normal"
$newPara.Shading.Texture = 0
$newPara.Shading.BackgroundPatternColor = 16382457
$newPara.Shading.ForegroundPatternColor = -16777216

Write-Host "Paragraph count: $($d.Paragraphs.Count)"
